# README updated with usage instructions
#
# The underlying data rows (2-5) get cycled: the former row 5 ("C", index 3)
# moves up to row 2, and the previous rows 2-4 ("A", "Aurobindo pharma", "B")
# shift down by one row into rows 3-5. Row 6 ("D") stays in place.
# In addition every data row (2-6) gets a freshly computed value in the
# "Debt equity ratio.1" (X) and "Score" (Y) columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2  (was row 5: "C", index 3) ----
$ws.Cells.Item(2, 1).Value = 3
$ws.Cells.Item(2, 2).Value = "C"
$ws.Cells.Item(2, 5).Value = 3
$ws.Cells.Item(2, 6).Value = 3
$ws.Cells.Item(2, 7).Value = ""
$ws.Cells.Item(2, 8).Value = 4
$ws.Cells.Item(2, 9).Value = 3
$ws.Cells.Item(2, 10).Value = 6
$ws.Cells.Item(2, 11).Value = 6
$ws.Cells.Item(2, 12).Value = 7
$ws.Cells.Item(2, 13).Value = 0
$ws.Cells.Item(2, 14).Value = 10
$ws.Cells.Item(2, 15).Value = 0
$ws.Cells.Item(2, 16).Value = 10
$ws.Cells.Item(2, 17).Value = 0
$ws.Cells.Item(2, 19).Value = 0
$ws.Cells.Item(2, 21).Value = 0
$ws.Cells.Item(2, 23).Value = 0
$ws.Cells.Item(2, 24).Value = 10
$ws.Cells.Item(2, 25).Value = 1110

# ---- Row 3  (was row 2: "A", index 0) ----
$ws.Cells.Item(3, 1).Value = 0
$ws.Cells.Item(3, 2).Value = "A"
$ws.Cells.Item(3, 5).Value = 7
$ws.Cells.Item(3, 6).Value = 2
$ws.Cells.Item(3, 7).Value = 12
$ws.Cells.Item(3, 8).Value = 6
$ws.Cells.Item(3, 9).Value = 8
$ws.Cells.Item(3, 10).Value = 7
$ws.Cells.Item(3, 11).Value = 8
$ws.Cells.Item(3, 12).Value = 5
$ws.Cells.Item(3, 13).Value = 111
$ws.Cells.Item(3, 14).Value = 1
$ws.Cells.Item(3, 15).Value = 13.7
$ws.Cells.Item(3, 16).Value = 5.620689655172414
$ws.Cells.Item(3, 17).Value = 1.43
$ws.Cells.Item(3, 19).Value = 19.591
$ws.Cells.Item(3, 21).Value = 0.54
$ws.Cells.Item(3, 23).Value = 0.33
$ws.Cells.Item(3, 24).Value = 7.444444444444444
$ws.Cells.Item(3, 25).Value = 783.1685823754789

# ---- Row 4  (was row 3: "Aurobindo pharma", index 1) ----
$ws.Cells.Item(4, 1).Value = 1
$ws.Cells.Item(4, 2).Value = "Aurobindo pharma"
$ws.Cells.Item(4, 5).Value = 4
$ws.Cells.Item(4, 6).Value = 1
$ws.Cells.Item(4, 8).Value = 4
$ws.Cells.Item(4, 9).Value = 7.5
$ws.Cells.Item(4, 10).Value = 7.5
$ws.Cells.Item(4, 11).Value = 7.5
$ws.Cells.Item(4, 12).Value = 7.5
$ws.Cells.Item(4, 13).Value = 13
$ws.Cells.Item(4, 14).Value = 9.666666666666666
$ws.Cells.Item(4, 15).Value = 9.01
$ws.Cells.Item(4, 16).Value = 7.237931034482759
$ws.Cells.Item(4, 17).Value = 1.53
$ws.Cells.Item(4, 19).Value = 13.7853
$ws.Cells.Item(4, 21).Value = 1.35
$ws.Cells.Item(4, 23).Value = 0.39
$ws.Cells.Item(4, 24).Value = 6.777777777777777
$ws.Cells.Item(4, 25).Value = 769.5386973180076

# ---- Row 5  (was row 4: "B", index 2) ----
$ws.Cells.Item(5, 1).Value = 2
$ws.Cells.Item(5, 2).Value = "B"
$ws.Cells.Item(5, 5).Value = 6
$ws.Cells.Item(5, 6).Value = 3
$ws.Cells.Item(5, 8).Value = 3
$ws.Cells.Item(5, 9).Value = 4
$ws.Cells.Item(5, 10).Value = 3.5
$ws.Cells.Item(5, 11).Value = 9
$ws.Cells.Item(5, 12).Value = 8
$ws.Cells.Item(5, 13).Value = 0
$ws.Cells.Item(5, 14).Value = 10
$ws.Cells.Item(5, 15).Value = 22
$ws.Cells.Item(5, 16).Value = 2.758620689655173
$ws.Cells.Item(5, 17).Value = 4.53
$ws.Cells.Item(5, 19).Value = 99.66
$ws.Cells.Item(5, 21).Value = 1.72
$ws.Cells.Item(5, 23).Value = 0.72
$ws.Cells.Item(5, 24).Value = 3.111111111111111
$ws.Cells.Item(5, 25).Value = 377.6628352490422

# ---- Row 6  ("D", index 4 — unchanged position, only X6/Y6 are refreshed) ----
$ws.Cells.Item(6, 1).Value = 4
$ws.Cells.Item(6, 2).Value = "D"
$ws.Cells.Item(6, 5).Value = 6
$ws.Cells.Item(6, 6).Value = 7
$ws.Cells.Item(6, 8).Value = 4
$ws.Cells.Item(6, 9).Value = 5
$ws.Cells.Item(6, 10).Value = 6
$ws.Cells.Item(6, 11).Value = 7
$ws.Cells.Item(6, 12).Value = 8
$ws.Cells.Item(6, 13).Value = 211
$ws.Cells.Item(6, 14).Value = 1
$ws.Cells.Item(6, 15).Value = 44
$ws.Cells.Item(6, 16).Value = 1
$ws.Cells.Item(6, 17).Value = 4.9
$ws.Cells.Item(6, 19).Value = 215.6
$ws.Cells.Item(6, 21).Value = 2.2
$ws.Cells.Item(6, 23).Value = 1.7
$ws.Cells.Item(6, 24).Value = 1
$ws.Cells.Item(6, 25).Value = 111
